$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 1-21, columns A-C. Column D (Score) is left blank
# except for style application, handled separately below.
$data = @(
    @("Team 1", "Team 2", "Winner"),
    @("Mora A", "Ruhuna", "Not Played Yet"),
    @("Kelani", "Pera", "Not Played Yet"),
    @("Colombo", "Wayamba", "Not Played Yet"),
    @("Japura", "Rajarata", "Not Played Yet"),
    @("Mora B", "Colombo", "Not Played Yet"),
    @("Wayamba", "Kelani", "Not Played Yet"),
    @("Ruhuna", "Japura", "Not Played Yet"),
    @("Colombo", "Pera", "Not Played Yet"),
    @("Mora B", "Wayamba", "Not Played Yet"),
    @("Mora A", "Rajarata", "Not Played Yet"),
    @("Colombo", "Kelani", "Not Played Yet"),
    @("Mora B", "Pera", "Not Played Yet"),
    @("Mora A", "Japura", "Not Played Yet"),
    @("Wayamba", "Pera", "Not Played Yet"),
    @("Mora B", "Kelani", "Not Played Yet"),
    @("Ruhuna", "Rajarata", "Not Played Yet"),
    @("WA 1st", "WB 2nd", "Not Played Yet"),
    @("WB 1st", "WA 2nd", "Not Played Yet"),
    @("37 Looser", "38 Looser", "Not Played Yet"),
    @("37 Winner", "38 Winner", "Not Played Yet")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 1
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}

# header row D1 keeps "Score" (unchanged)
$ws.Cells.Item(1, 4).Value = "Score"

# Clear any leftover D values below the header, and ensure D3:D11 carry
# the numeric-text style (style index 2 in styles.xml) with no value.
for ($r = 2; $r -le 21; $r++) {
    $ws.Cells.Item($r, 4).Value = $null
}
$ws.Range("D3:D11").NumberFormat = "@"

# Column widths (A=11.21875, B=11.88671875, C=14.109375 characters)
$ws.Columns.Item(1).ColumnWidth = 11.21875
$ws.Columns.Item(2).ColumnWidth = 11.88671875
$ws.Columns.Item(3).ColumnWidth = 14.109375

# Selection matches the post-edit state
$ws.Range("G9").Select()
